$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.931.27'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '2.603.11'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '''522.73'
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').Value = '''154.82'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +2.39%  '
$ws.Range('E9').Value = '  +2.20%  '
$ws.Range('E10').Value = '  +2.03%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '3.059.47'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').Value = '60.897.85'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D17').Value = '2.604.29'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '''354.63'
$ws.Range('E19').Value = '  +3.21%  '
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('D23').Value = '''60.94'
$ws.Range('E23').Value = '  +2.07%  '
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('E25').Value = '  +0.71%  '
$ws.Range('D26').Value = '2.718.71'
$ws.Range('E26').Value = '  +1.28%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '''6.28'
$ws.Range('E31').Value = '  +10.39%  '
$ws.Range('D32').Value = '''19.41'
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('D34').Value = '''148.05'
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('D35').Value = '''4.19'
$ws.Range('E35').Value = '  +5.42%  '
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('D37').Value = '''0.917'
$ws.Range('E37').Value = '  +8.59%  '
$ws.Range('D38').Value = '''0.876'
$ws.Range('E38').Value = '  +3.66%  '
$ws.Range('E39').Value = '  +2.78%  '
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '''289.91'
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('E43').Value = '  +2.52%  '
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '''19.59'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('D50').Value = '''10.34'
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '''19.20'
$ws.Range('E51').Value = '  +9.33%  '
